# Append the new EUR->ARS quote row (2025-10-13 15:20:30) to the bottom of the
# rate-history sheet, one row below the previous last entry (row 74 -> row 75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Column A holds an ISO "yyyy-mm-dd" string. Excel's automatic type inference
# would otherwise turn that literal text into a date serial, so force it to
# stay plain text with a leading quote-prefix (the standard Excel idiom),
# exactly like the rest of the historical rows in this sheet.
$ws.Cells.Item($row, 1).Value = "'2025-10-13"
$ws.Cells.Item($row, 2).Value = "15:20:30"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,753.3027"
